$d = $word.ActiveDocument
$d.Content.Find.Execute("antennas", $true, $false, $false, $false, $false, $true, 1, $false, "transmissions", 2)
